$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.789.23"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "3.176.10"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'594.12"
$ws.Range("E5").Value = "  +3.48%  "
$ws.Range("D6").Value = "'152.48"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.175.40"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").Value = "'6.05"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "'0.511"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "'38.68"
$ws.Range("E14").Value = "  +3.93%  "
$ws.Range("D15").Value = "3.695.59"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "65.860.73"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "'7.40"
$ws.Range("E17").Value = "  +4.10%  "
$ws.Range("D18").Value = "3.172.86"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "'506.27"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "'15.30"
$ws.Range("E21").Value = "  +2.56%  "
$ws.Range("D22").Value = "'0.732"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("D23").Value = "'7.97"
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("D24").Value = "'14.91"
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("D25").Value = "'84.55"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'9.21"
$ws.Range("E27").Value = "  +3.30%  "
$ws.Range("D28").Value = "'2.97"
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("D29").Value = "'2.28"
$ws.Range("E29").Value = "  +4.83%  "
$ws.Range("D30").Value = "'6.96"
$ws.Range("E30").Value = "  +12.58%  "
$ws.Range("D31").Value = "'2.88"
$ws.Range("E31").Value = "  +2.84%  "
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "'6.45"
$ws.Range("E35").Value = "  -1.19%  "
$ws.Range("D36").Value = "'54.67"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "'0.0899"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "'480.07"
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("D39").Value = "'0.0417"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").Value = "'8.79"
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("D41").Value = "'2.84"
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("D43").Value = "'0.298"
$ws.Range("E43").Value = "  +5.25%  "
$ws.Range("D44").Value = "0.0₃0649"
$ws.Range("E44").Value = "  +10.42%  "
$ws.Range("D45").Value = "2.882.68"
$ws.Range("E45").Value = "  -5.78%  "
$ws.Range("D46").Value = "'2.40"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("D47").Value = "'28.26"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").Value = "'2.30"
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("E51").Value = "  +5.47%  "
